# Updated capital structure database
# Applies the updated row 2 (summary/aggregate row) and row 3 (company row)
# values for the Russia / Brokerage & Investment Banking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("D2").Value  = 0.139
$ws.Range("E2").Value  = 0.0314

$ws.Range("K2").Value  = 0.801
$ws.Range("L2").Value  = 0.1849884526558891
$ws.Range("M2").Value  = 0
$ws.Range("N2").Value  = 0
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 0
$ws.Range("Q2").Value  = 0
$ws.Range("R2").Value  = 0

# Column T is dropped entirely for this row (removed from the sheet data)
$ws.Range("T2").ClearContents()

$ws.Range("U2").Value  = 20.2
$ws.Range("V2").Value  = 0.464367816091954
$ws.Range("W2").Value  = 0.04914110429447853
$ws.Range("X2").Value  = 0.03972798311003502
$ws.Range("Y2").Value  = 0.009413121184443513
$ws.Range("Z2").Value  = 2.152087475149105

$ws.Range("AB2").Value = 0.03946404969782924
$ws.Range("AC2").Value = -0.03946404969782924
$ws.Range("AD2").Value = 3.52
$ws.Range("AF2").Value = 3.52
$ws.Range("AG2").Value = -16.68
$ws.Range("AH2").Value = 0.07486176095278604
$ws.Range("AI2").Value = 0.1964285714285714
$ws.Range("AJ2").Value = -0.6219239373601789
$ws.Range("AK2").Value = 7.315789473684212

# ---- Row 3 ----
$ws.Range("D3").Value  = 0.139
$ws.Range("E3").Value  = 0.0314

$ws.Range("K3").Value  = 0.801
$ws.Range("L3").Value  = 0.1849884526558891
$ws.Range("M3").Value  = -0
$ws.Range("N3").Value  = -0
$ws.Range("O3").Value  = -0
$ws.Range("P3").Value  = -0
$ws.Range("Q3").Value  = -0
$ws.Range("R3").Value  = -0

# Column T is dropped entirely for this row (removed from the sheet data)
$ws.Range("T3").ClearContents()

$ws.Range("U3").Value  = 20.2
$ws.Range("V3").Value  = 0.464367816091954
$ws.Range("W3").Value  = 0.04914110429447853
$ws.Range("X3").Value  = 0.03972798311003502
$ws.Range("Y3").Value  = 0.009413121184443513
$ws.Range("Z3").Value  = 2.152087475149105

$ws.Range("AB3").Value = 0.03946404969782924
$ws.Range("AC3").Value = -0.03946404969782924
$ws.Range("AD3").Value = 3.52
$ws.Range("AF3").Value = 3.52
$ws.Range("AG3").Value = -16.68
$ws.Range("AH3").Value = 0.07486176095278604
$ws.Range("AI3").Value = 0.1964285714285714
$ws.Range("AJ3").Value = -0.6219239373601789
$ws.Range("AK3").Value = 7.315789473684212
